$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Update the numeric table values (the BFGS re-fit changed the reported
#    ratios/CIs for several cells; two cells -- 0.65 (0.64 - 0.65) and
#    0.3 (0.3 - 0.3) -- are unchanged by the new optimizer and are left as-is).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("2.8 (0.63 - 14.9)", $true, $false, $false, $false, $false, `
    $true, 1, $false, "1.43 (0.76 - 2.37)", 2)
$d.Content.Find.Execute("0.2 (0.01 - 0.62)", $true, $false, $false, $false, $false, `
    $true, 1, $false, "0.17 (0.05 - 0.33)", 2)
$d.Content.Find.Execute("0.89 (0.22 - 4.79)", $true, $false, $false, $false, $false, `
    $true, 1, $false, "0.34 (0.25 - 0.64)", 2)
$d.Content.Find.Execute("2.26 (0.5 - 33.64)", $true, $false, $false, $false, $false, `
    $true, 1, $false, "0.7 (0.6 - 0.88)", 2)
$d.Content.Find.Execute("0.05 (0.01 - 0.16)", $true, $false, $false, $false, $false, `
    $true, 1, $false, "0.07 (0.04 - 0.09)", 2)
$d.Content.Find.Execute("2.05 (0.45 - 6.89)", $true, $false, $false, $false, $false, `
    $true, 1, $false, "0.72 (0.56 - 1.02)", 2)

# ---------------------------------------------------------------------------
# 2. Switch the table's font from Arial to Helvetica everywhere.
#    a) A formatting-only Find/Replace rewrites every run that currently uses
#       Arial so that ascii/hAnsi/eastAsia/cs all become Helvetica.
# ---------------------------------------------------------------------------
$find = $d.Content.Find
$find.ClearFormatting()
$find.Font.Name = "Arial"
$find.Replacement.ClearFormatting()
$find.Replacement.Font.Name = "Helvetica"
$find.Replacement.Font.NameFarEast = "Helvetica"
$find.Replacement.Font.NameOther = "Helvetica"
$find.Replacement.Font.NameBi = "Helvetica"
$find.Text = ""
$find.Replacement.Text = ""
$find.Forward = $true
$find.Wrap = 1
$find.Format = $true
$find.MatchWildcards = $false
$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $true, $null, 2)

#    b) Sweep every table's full range (header + body, including paragraph
#       marks) so any remaining Arial references pick up the new font too.
foreach ($t in $d.Tables) {
    $t.Range.Font.Name = "Helvetica"
}

#    c) Belt-and-braces: walk every paragraph in the document as well.
foreach ($p in $d.Paragraphs) {
    $p.Range.Font.Name = "Helvetica"
}
